# Remove Sheet2 entirely
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Delete()

# Update the remaining Sheet1 content: translate headers and adjust trailing
# whitespace on the String10..String19 labels.
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").Value = "StringTitle"
$ws1.Range("B1").Value = "DateTitle"
$ws1.Range("C1").Value = "DoubleTitle"

$ws1.Range("A2").Value  = "String10 "
$ws1.Range("A3").Value  = "String11 "
$ws1.Range("A4").Value  = "String12 "
$ws1.Range("A5").Value  = "String13 "
$ws1.Range("A6").Value  = "String14 "

$ideo = [char]0x3000
$ws1.Range("A7").Value  = "String15" + $ideo
$ws1.Range("A8").Value  = "String16" + $ideo
$ws1.Range("A9").Value  = "String17" + $ideo

$ws1.Range("A10").Value = "String18 "
$ws1.Range("A11").Value = "String19 "

# Update the view state: Sheet1 becomes the selected/active tab with a new
# selected cell, matching what the workbook looked like after trimming it
# down to a single sheet.
[void]$ws1.Range("B13").Select()

Write-Host "done"
